$d = $word.ActiveDocument

$replacements = @(
    @{old = "2025-10-20 Monday"; new = "2025-10-21 Tuesday"},
    @{old = "523×5="; new = "315×7="},
    @{old = "316×8="; new = "455×9="},
    @{old = "534×9="; new = "510×9="},
    @{old = "800×8="; new = "320×5="},
    @{old = "554×7="; new = "920×3="},
    @{old = "518×2="; new = "564×7="},
    @{old = "875×5="; new = "822×9="},
    @{old = "571×5="; new = "272×4="},
    @{old = "951×6="; new = "321×6="},
    @{old = "478×5="; new = "361×3="},
    @{old = "685×5="; new = "579×6="},
    @{old = "442×7="; new = "218×8="},
    @{old = "887×4="; new = "363×3="},
    @{old = "387×2="; new = "511×4="},
    @{old = "206×5="; new = "563×6="},
    @{old = "522×4="; new = "845×3="},
    @{old = "442×4="; new = "227×9="},
    @{old = "525×6="; new = "501×5="},
    @{old = "726×7="; new = "612×6="},
    @{old = "601×3="; new = "320×5="},
    @{old = "219×5="; new = "848×9="},
    @{old = "297×7="; new = "110×2="},
    @{old = "414×8="; new = "795×7="},
    @{old = "102×3="; new = "338×4="},
    @{old = "898×5="; new = "795×5="}
)

foreach ($r in $replacements) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Replacement.ClearFormatting()
    $find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
